$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.705.27'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.303.54'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.45%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.83'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.04'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.32%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.606'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.85'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.16%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.07%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.63%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.994'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.18%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.39'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.74%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.650.80'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.306.99'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.616.59'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.97%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.62'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.85%  '

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000106'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.84%  '

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.75'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +31.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.04'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.02%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.77'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.61%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.01%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.91'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.83%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.20%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.65'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.60'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +13.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.76'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.95%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '165.73'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.37%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.56%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.65%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.63%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.114'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.59'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.45%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.81%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.74'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.27%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.26%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +11.83%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.53'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '96.14'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.18%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.63%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '117.94'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.50%  '

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.45'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '80.56'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.74%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.651.79'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.89'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.18%  '
